$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "penismaster"
$ws.Range("B2").Value = "5fa285e1bebe0a6623e33afc04a1fbd5"

$ws.Range("A3").Value = 123
$ws.Range("B3").Value = "5fa285e1bebe0a6623e33afc04a1fbd5"

$ws.Range("A4").Value = "wdf"
$ws.Range("B4").Value = "d1c364f4b712dc6c804c2b773de76d97"

$ws.Range("A5").Value = "penis"
$ws.Range("B5").Value = "46dc363a5e754c6781f8889094b288c4"
